# The codeforiati:category-name / codeforiati:group-name columns (D/E) were
# swapped, and the codeforiati:group-code / codeforiati:category-code
# columns (F/G) were swapped, for every row of the sheet (including the
# header row). Use Range.Copy so that text-typed cells (e.g. the numeric
# looking "110", "111" codes) stay stored as text instead of being
# reinterpreted as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

$colD = "D" + $firstRow + ":D" + $lastRow
$colE = "E" + $firstRow + ":E" + $lastRow
$colF = "F" + $firstRow + ":F" + $lastRow
$colG = "G" + $firstRow + ":G" + $lastRow

# Scratch columns well outside the used range of the sheet.
$scratch1 = "ZZ" + $firstRow + ":ZZ" + $lastRow
$scratch2 = "ZY" + $firstRow + ":ZY" + $lastRow

# Swap D <-> E using a scratch column as temporary holder.
$ws.Range($colD).Copy($ws.Range($scratch1))
$ws.Range($colE).Copy($ws.Range($colD))
$ws.Range($scratch1).Copy($ws.Range($colE))

# Swap F <-> G using a scratch column as temporary holder.
$ws.Range($colF).Copy($ws.Range($scratch2))
$ws.Range($colG).Copy($ws.Range($colF))
$ws.Range($scratch2).Copy($ws.Range($colG))

# Clean up the scratch columns so they don't linger in the saved file.
$ws.Range($scratch1).Clear()
$ws.Range($scratch2).Clear()
